$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Discount Rate (r)" column (C) entries that were placeholder
# percentage-range strings with actual computed numeric percentages,
# formatted as percentages (reusing the existing 0.00% style already used
# by C8).
$ws.Range("C5").Value = 0.1004
$ws.Range("C5").NumberFormat = "0.00%"

$ws.Range("C6").Value = 0.1054
$ws.Range("C6").NumberFormat = "0.00%"

$ws.Range("C7").Value = 0.1014
$ws.Range("C7").NumberFormat = "0.00%"

$ws.Range("C8").Value = 0.102
$ws.Range("C8").NumberFormat = "0.00%"

$ws.Range("C9").Value = 0.1004
$ws.Range("C9").NumberFormat = "0.00%"

$ws.Range("C10").Value = 0.1004
$ws.Range("C10").NumberFormat = "0.00%"

$ws.Range("C11").Value = 0.1004
$ws.Range("C11").NumberFormat = "0.00%"

$ws.Range("C12").Value = 0.1004
$ws.Range("C12").NumberFormat = "0.00%"

$ws.Range("C13").Value = 0.1004
$ws.Range("C13").NumberFormat = "0.00%"

$ws.Range("C14").Value = 0.1004
$ws.Range("C14").NumberFormat = "0.00%"

# Move the current selection / view so the sheet no longer shows the
# scrolled-down, E13-selected state.
$ws.Range("A17").Select()
